$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '25.979.50'
$ws.Cells.Item(2, 5).Value = '  -0.93%  '
$ws.Cells.Item(3, 4).Value = '1.745.82'
$ws.Cells.Item(3, 5).Value = '  -0.25%  '
$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = '1.000'
$ws.Cells.Item(4, 4).NumberFormat = "General"
$ws.Cells.Item(4, 5).Value = '  +0.09%  '
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '248.72'
$ws.Cells.Item(5, 4).NumberFormat = "General"
$ws.Cells.Item(5, 5).Value = '  +2.34%  '
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '0.9998'
$ws.Cells.Item(6, 4).NumberFormat = "General"
$ws.Cells.Item(6, 5).Value = '  -0.02%  '
$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.5084'
$ws.Cells.Item(7, 4).NumberFormat = "General"
$ws.Cells.Item(7, 5).Value = '  -8.38%  '
$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.2765'
$ws.Cells.Item(8, 4).NumberFormat = "General"
$ws.Cells.Item(8, 5).Value = '  -2.87%  '
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.06192'
$ws.Cells.Item(9, 4).NumberFormat = "General"
$ws.Cells.Item(9, 5).Value = '  -0.12%  '
$ws.Cells.Item(10, 4).Value = '1.747.81'
$ws.Cells.Item(10, 5).Value = '  -0.22%  '
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '0.07268'
$ws.Cells.Item(11, 4).NumberFormat = "General"
$ws.Cells.Item(11, 5).Value = '  +0.99%  '
$ws.Cells.Item(12, 2).Value = 'Solana'
$ws.Cells.Item(12, 3).Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '15.21'
$ws.Cells.Item(12, 4).NumberFormat = "General"
$ws.Cells.Item(12, 5).Value = '  -2.85%  '
$ws.Cells.Item(13, 2).Value = 'Polygon'
$ws.Cells.Item(13, 3).Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '0.6541'
$ws.Cells.Item(13, 4).NumberFormat = "General"
$ws.Cells.Item(13, 5).Value = '  -1.59%  '
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '4.684'
$ws.Cells.Item(14, 4).NumberFormat = "General"
$ws.Cells.Item(14, 5).Value = '  +0.38%  '
$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '77.74'
$ws.Cells.Item(15, 4).NumberFormat = "General"
$ws.Cells.Item(15, 5).Value = '  -1.19%  '
$ws.Cells.Item(16, 5).Value = '  +0.07%  '
$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.9999'
$ws.Cells.Item(17, 4).NumberFormat = "General"
$ws.Cells.Item(17, 5).Value = '  +0.07%  '
$ws.Cells.Item(18, 4).Value = '25.989.33'
$ws.Cells.Item(18, 5).Value = '  -0.57%  '
$ws.Cells.Item(19, 5).Value = '  -0.18%  '
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '0.000006885'
$ws.Cells.Item(20, 4).NumberFormat = "General"
$ws.Cells.Item(20, 5).Value = '  +0.86%  '
$ws.Cells.Item(21, 4).Value = '1.970.92'
$ws.Cells.Item(21, 5).Value = '  -0.61%  '
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '4.481'
$ws.Cells.Item(22, 4).NumberFormat = "General"
$ws.Cells.Item(22, 5).Value = '  +1.39%  '
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '8.761'
$ws.Cells.Item(23, 4).NumberFormat = "General"
$ws.Cells.Item(23, 5).Value = '  -0.62%  '
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '5.386'
$ws.Cells.Item(24, 4).NumberFormat = "General"
$ws.Cells.Item(24, 5).Value = '  +1.03%  '
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '136.36'
$ws.Cells.Item(25, 4).NumberFormat = "General"
$ws.Cells.Item(25, 5).Value = '  -2.85%  '
$ws.Cells.Item(26, 5).Value = '  -0.24%  '
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '15.30'
$ws.Cells.Item(27, 4).NumberFormat = "General"
$ws.Cells.Item(27, 5).Value = '  -0.54%  '
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = '1.787'
$ws.Cells.Item(28, 4).NumberFormat = "General"
$ws.Cells.Item(28, 5).Value = '  -2.02%  '
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '105.82'
$ws.Cells.Item(29, 4).NumberFormat = "General"
$ws.Cells.Item(29, 5).Value = '  -0.26%  '
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = '3.877'
$ws.Cells.Item(30, 4).NumberFormat = "General"
$ws.Cells.Item(30, 5).Value = '  +1.82%  '
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '0.08221'
$ws.Cells.Item(31, 4).NumberFormat = "General"
$ws.Cells.Item(31, 5).Value = '  -3.82%  '
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '3.654'
$ws.Cells.Item(32, 4).NumberFormat = "General"
$ws.Cells.Item(32, 5).Value = '  -0.13%  '
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.04664'
$ws.Cells.Item(33, 4).NumberFormat = "General"
$ws.Cells.Item(33, 5).Value = '  +0.21%  '
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '2.653'
$ws.Cells.Item(34, 4).NumberFormat = "General"
$ws.Cells.Item(34, 5).Value = '  -0.26%  '
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '0.9984'
$ws.Cells.Item(35, 4).NumberFormat = "General"
$ws.Cells.Item(35, 5).Value = '  -1.37%  '
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '0.6145'
$ws.Cells.Item(36, 4).NumberFormat = "General"
$ws.Cells.Item(36, 5).Value = '  -2.64%  '
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '2.781'
$ws.Cells.Item(37, 4).NumberFormat = "General"
$ws.Cells.Item(37, 5).Value = '  +2.78%  '
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '0.01618'
$ws.Cells.Item(38, 4).NumberFormat = "General"
$ws.Cells.Item(38, 5).Value = '  +0.07%  '
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '1.936'
$ws.Cells.Item(39, 4).NumberFormat = "General"
$ws.Cells.Item(39, 5).Value = '  -2.41%  '
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '0.9995'
$ws.Cells.Item(40, 4).NumberFormat = "General"
$ws.Cells.Item(40, 5).Value = '  -0.06%  '
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '100.58'
$ws.Cells.Item(41, 4).NumberFormat = "General"
$ws.Cells.Item(41, 5).Value = '  +0.94%  '
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '0.3932'
$ws.Cells.Item(42, 4).NumberFormat = "General"
$ws.Cells.Item(42, 5).Value = '  -0.83%  '
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '0.7680'
$ws.Cells.Item(43, 4).NumberFormat = "General"
$ws.Cells.Item(43, 5).Value = '  +1.99%  '
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '5.018'
$ws.Cells.Item(44, 4).NumberFormat = "General"
$ws.Cells.Item(44, 5).Value = '  -0.58%  '
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.1156'
$ws.Cells.Item(45, 4).NumberFormat = "General"
$ws.Cells.Item(45, 5).Value = '  +0.11%  '
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '6.358'
$ws.Cells.Item(46, 4).NumberFormat = "General"
$ws.Cells.Item(46, 5).Value = '  -0.58%  '
$ws.Cells.Item(47, 2).Value = 'Cronos'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.05344'
$ws.Cells.Item(47, 4).NumberFormat = "General"
$ws.Cells.Item(47, 5).Value = '  +0.04%  '
$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '55.68'
$ws.Cells.Item(48, 4).NumberFormat = "General"
$ws.Cells.Item(48, 5).Value = '  +1.05%  '
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '30.70'
$ws.Cells.Item(49, 4).NumberFormat = "General"
$ws.Cells.Item(49, 5).Value = '  -1.11%  '
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '7.614'
$ws.Cells.Item(50, 4).NumberFormat = "General"
$ws.Cells.Item(50, 5).Value = '  -0.66%  '
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '0.3453'
$ws.Cells.Item(51, 4).NumberFormat = "General"
$ws.Cells.Item(51, 5).Value = '  -1.94%  '
